$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XRAY")

# Row 4 (Inventory): update B4:F4
$ws.Range("B4").Value = 466000000.0
$ws.Range("C4").Value = 489000000.0
$ws.Range("D4").Value = 549000000.0
$ws.Range("E4").Value = 591000000.0
$ws.Range("F4").Value = 562000000.0

# Row 14 (Accounts Payable): update B14:F14
$ws.Range("B14").Value = 305000000.0
$ws.Range("C14").Value = 243000000.0
$ws.Range("D14").Value = 214000000.0
$ws.Range("E14").Value = 271000000.0
$ws.Range("F14").Value = 308000000.0

# Row 21 (Long Term Tax Liability (Deferred)): update B21:F21
$ws.Range("B21").Value = 384000000.0
$ws.Range("C21").Value = 423000000.0
$ws.Range("D21").Value = 439000000.0
$ws.Range("E21").Value = 465000000.0
$ws.Range("F21").Value = 467000000.0
